$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Formula = '="87548754"'
$ws.Range("B2").Formula = '="yanet altamirano quiroz"'
$ws.Range("C2").Formula = '="2025-03-15"'
$ws.Range("D2").Formula = '="07:50:27"'
$ws.Range("E2").Formula = '="13:05:18"'
$ws.Range("F2").Formula = '="14:03:07"'
$ws.Range("G2").Formula = '="18:44:33"'
$ws.Range("H2").Formula = '="0 minutos"'
$ws.Range("A2:H2").Copy()
$ws.Range("A2:H2").PasteSpecial(-4163)

# Row 3
$ws.Range("A3").Formula = '="87548754"'
$ws.Range("B3").Formula = '="yanet altamirano quiroz"'
$ws.Range("C3").Formula = '="2025-03-18"'
$ws.Range("D3").Formula = '="07:53:38"'
$ws.Range("E3").Formula = '="13:53:38"'
$ws.Range("A3:E3").Copy()
$ws.Range("A3:E3").PasteSpecial(-4163)
$ws.Range("G3").Formula = '="17:53:39"'
$ws.Range("H3").Formula = '="0 minutos"'
$ws.Range("I3").Formula = '="Tiene horas sin marcar"'
$ws.Range("G3:I3").Copy()
$ws.Range("G3:I3").PasteSpecial(-4163)

# Row 4
$ws.Range("A4").Formula = '="87548754"'
$ws.Range("B4").Formula = '="yanet altamirano quiroz"'
$ws.Range("C4").Formula = '="2025-03-17"'
$ws.Range("D4").Formula = '="07:52:27"'
$ws.Range("A4:D4").Copy()
$ws.Range("A4:D4").PasteSpecial(-4163)
$ws.Range("F4").Formula = '="14:14:52"'
$ws.Range("G4").Formula = '="18:10:06"'
$ws.Range("H4").Formula = '="0 minutos"'
$ws.Range("I4").Formula = '="Tiene horas sin marcar"'
$ws.Range("F4:I4").Copy()
$ws.Range("F4:I4").PasteSpecial(-4163)

# Row 5
$ws.Range("A5").Formula = '="87548754"'
$ws.Range("B5").Formula = '="yanet altamirano quiroz"'
$ws.Range("C5").Formula = '="2025-03-16"'
$ws.Range("D5").Formula = '="07:51:52"'
$ws.Range("E5").Formula = '="13:07:12"'
$ws.Range("F5").Formula = '="14:06:03"'
$ws.Range("G5").Formula = '="17:53:26"'
$ws.Range("H5").Formula = '="0 minutos"'
$ws.Range("A5:H5").Copy()
$ws.Range("A5:H5").PasteSpecial(-4163)

# Row 6
$ws.Range("A6").Formula = '="12345667"'
$ws.Range("B6").Formula = '="carla siares adrianzen"'
$ws.Range("C6").Formula = '="2025-03-16"'
$ws.Range("D6").Formula = '="07:51:52"'
$ws.Range("E6").Formula = '="13:07:12"'
$ws.Range("F6").Formula = '="14:06:03"'
$ws.Range("G6").Formula = '="17:53:26"'
$ws.Range("H6").Formula = '="0 minutos"'
$ws.Range("A6:H6").Copy()
$ws.Range("A6:H6").PasteSpecial(-4163)

# Row 7
$ws.Range("A7").Formula = '="12345667"'
$ws.Range("B7").Formula = '="carla siares adrianzen"'
$ws.Range("C7").Formula = '="2025-03-17"'
$ws.Range("D7").Formula = '="07:52:27"'
$ws.Range("A7:D7").Copy()
$ws.Range("A7:D7").PasteSpecial(-4163)
$ws.Range("F7").Formula = '="14:14:52"'
$ws.Range("G7").Formula = '="18:10:06"'
$ws.Range("H7").Formula = '="0 minutos"'
$ws.Range("I7").Formula = '="Tiene horas sin marcar"'
$ws.Range("F7:I7").Copy()
$ws.Range("F7:I7").PasteSpecial(-4163)

# Row 8
$ws.Range("A8").Formula = '="12345667"'
$ws.Range("B8").Formula = '="carla siares adrianzen"'
$ws.Range("C8").Formula = '="2025-03-18"'
$ws.Range("D8").Formula = '="07:53:38"'
$ws.Range("A8:D8").Copy()
$ws.Range("A8:D8").PasteSpecial(-4163)
$ws.Range("H8").Formula = '="0 minutos"'
$ws.Range("I8").Formula = '="Tiene horas sin marcar"'
$ws.Range("H8:I8").Copy()
$ws.Range("H8:I8").PasteSpecial(-4163)

# Row 9
$ws.Range("A9").Formula = '="12345667"'
$ws.Range("B9").Formula = '="carla siares adrianzen"'
$ws.Range("C9").Formula = '="2025-03-15"'
$ws.Range("D9").Formula = '="07:50:27"'
$ws.Range("E9").Formula = '="13:05:18"'
$ws.Range("F9").Formula = '="14:03:07"'
$ws.Range("G9").Formula = '="18:44:33"'
$ws.Range("H9").Formula = '="0 minutos"'
$ws.Range("A9:H9").Copy()
$ws.Range("A9:H9").PasteSpecial(-4163)

# Clear blank cells (no value)
$ws.Range("F3").ClearContents()
$ws.Range("E4").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("G8").ClearContents()

# Empty-string cells (t="s" pointing at blank shared string)
$ws.Range("I2").Value = "'"
$ws.Range("I2").ClearFormats()
$ws.Range("I5").Value = "'"
$ws.Range("I5").ClearFormats()
$ws.Range("I6").Value = "'"
$ws.Range("I6").ClearFormats()
$ws.Range("I9").Value = "'"
$ws.Range("I9").ClearFormats()

$excel.CutCopyMode = $false